$d = $word.ActiveDocument

$pairs = @(
    @("140÷5=28, 0", "666÷6=111, 0"),
    @("628÷4=157, 0", "543÷6=90, 3"),
    @("987÷3=329, 0", "301÷8=37, 5"),
    @("812÷9=90, 2", "549÷9=61, 0"),
    @("478÷3=159, 1", "506÷7=72, 2"),
    @("123÷8=15, 3", "953÷8=119, 1"),
    @("264÷8=33, 0", "672÷7=96, 0"),
    @("524÷6=87, 2", "962÷5=192, 2"),
    @("218÷5=43, 3", "307÷5=61, 2"),
    @("434÷9=48, 2", "858÷3=286, 0"),
    @("955÷5=191, 0", "591÷8=73, 7"),
    @("653÷9=72, 5", "921÷2=460, 1"),
    @("918÷5=183, 3", "656÷7=93, 5"),
    @("428÷8=53, 4", "711÷5=142, 1"),
    @("290÷2=145, 0", "727÷7=103, 6"),
    @("686÷7=98, 0", "595÷8=74, 3"),
    @("904÷4=226, 0", "174÷9=19, 3"),
    @("883÷2=441, 1", "879÷5=175, 4"),
    @("509÷8=63, 5", "598÷6=99, 4"),
    @("607÷4=151, 3", "645÷2=322, 1"),
    @("890÷4=222, 2", "625÷4=156, 1"),
    @("873÷2=436, 1", "273÷6=45, 3"),
    @("840÷5=168, 0", "578÷8=72, 2"),
    @("291÷8=36, 3", "311÷3=103, 2"),
    @("565÷6=94, 1", "828÷8=103, 4")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
